$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

Set-TextValue $ws.Range('D2') '57.298.23'
Set-TextValue $ws.Range('E2') '  -1.24%  '
Set-TextValue $ws.Range('D3') '3.098.17'
Set-TextValue $ws.Range('E3') '  +0.74%  '
Set-TextValue $ws.Range('E4') '  -0.03%  '
Set-TextValue $ws.Range('D5') '523.04'
Set-TextValue $ws.Range('E5') '  +0.93%  '
Set-TextValue $ws.Range('D6') '140.77'
Set-TextValue $ws.Range('E6') '  -1.60%  '
Set-TextValue $ws.Range('E7') '  -0.08%  '
Set-TextValue $ws.Range('D8') '3.095.39'
Set-TextValue $ws.Range('E8') '  +0.73%  '
Set-TextValue $ws.Range('E9') '  +0.07%  '
Set-TextValue $ws.Range('D10') '7.19'
Set-TextValue $ws.Range('E10') '  -1.61%  '
Set-TextValue $ws.Range('D11') '0.108'
Set-TextValue $ws.Range('E11') '  -0.15%  '
Set-TextValue $ws.Range('E12') '  +1.52%  '
Set-TextValue $ws.Range('D13') '3.628.25'
Set-TextValue $ws.Range('E13') '  +0.79%  '
Set-TextValue $ws.Range('E14') '  +1.03%  '
Set-TextValue $ws.Range('D15') '26.01'
Set-TextValue $ws.Range('E15') '  -1.05%  '
Set-TextValue $ws.Range('E16') '  -0.94%  '
Set-TextValue $ws.Range('D17') '57.368.11'
Set-TextValue $ws.Range('E17') '  -1.14%  '
Set-TextValue $ws.Range('D18') '3.095.40'
Set-TextValue $ws.Range('E18') '  +0.62%  '
Set-TextValue $ws.Range('D19') '6.11'
Set-TextValue $ws.Range('E19') '  -0.16%  '
Set-TextValue $ws.Range('D20') '12.80'
Set-TextValue $ws.Range('E20') '  -0.96%  '
Set-TextValue $ws.Range('D21') '8.05'
Set-TextValue $ws.Range('E21') '  -0.95%  '
Set-TextValue $ws.Range('D22') '337.49'
Set-TextValue $ws.Range('E22') '  +1.12%  '
Set-TextValue $ws.Range('E23') '  -0.11%  '
Set-TextValue $ws.Range('D24') '0.513'
Set-TextValue $ws.Range('E24') '  +2.25%  '
Set-TextValue $ws.Range('D25') '66.62'
Set-TextValue $ws.Range('E25') '  +1.41%  '
Set-TextValue $ws.Range('E26') '  -1.51%  '
Set-TextValue $ws.Range('E27') '  +0.28%  '
Set-TextValue $ws.Range('D28') '0.0₃0910'
Set-TextValue $ws.Range('E28') '  +0.01%  '
Set-TextValue $ws.Range('D29') '6.49'
Set-TextValue $ws.Range('E29') '  -0.19%  '
Set-TextValue $ws.Range('D31') '7.17'
Set-TextValue $ws.Range('E31') '  -1.55%  '
Set-TextValue $ws.Range('E32') '  +2.13%  '
Set-TextValue $ws.Range('D33') '20.94'
Set-TextValue $ws.Range('E33') '  +0.80%  '
Set-TextValue $ws.Range('D34') '1.19'
Set-TextValue $ws.Range('E34') '  -0.92%  '
Set-TextValue $ws.Range('D35') '156.86'
Set-TextValue $ws.Range('E35') '  +1.27%  '
Set-TextValue $ws.Range('D36') '4.62'
Set-TextValue $ws.Range('E36') '  +0.98%  '
Set-TextValue $ws.Range('D37') '6.10'
Set-TextValue $ws.Range('E37') '  +1.34%  '
Set-TextValue $ws.Range('D38') '27.07'
Set-TextValue $ws.Range('E38') '  +0.13%  '
Set-TextValue $ws.Range('E39') '  -0.35%  '
Set-TextValue $ws.Range('D40') '0.0657'
Set-TextValue $ws.Range('E40') '  -3.10%  '
Set-TextValue $ws.Range('B41') 'Filecoin'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D41') '3.95'
Set-TextValue $ws.Range('E41') '  +0.65%  '
Set-TextValue $ws.Range('B42') 'RenzoRestakedETH'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue $ws.Range('D42') '3.138.58'
Set-TextValue $ws.Range('E42') '  +0.79%  '
Set-TextValue $ws.Range('E43') '  +3.94%  '
Set-TextValue $ws.Range('E44') '  +10.00%  '
Set-TextValue $ws.Range('D45') '36.60'
Set-TextValue $ws.Range('E45') '  +0.31%  '
Set-TextValue $ws.Range('D46') '0.999'
Set-TextValue $ws.Range('E46') '  -0.07%  '
Set-TextValue $ws.Range('D47') '2.308.64'
Set-TextValue $ws.Range('E47') '  +1.72%  '
Set-TextValue $ws.Range('D48') '0.0259'
Set-TextValue $ws.Range('E48') '  -0.50%  '
Set-TextValue $ws.Range('D49') '0.969'
Set-TextValue $ws.Range('E49') '  +2.25%  '
Set-TextValue $ws.Range('D50') '20.67'
Set-TextValue $ws.Range('E50') '  -2.05%  '
Set-TextValue $ws.Range('D51') '6.01'
Set-TextValue $ws.Range('E51') '  +1.22%  '
